# ArrayFormulaFunctions.xlsx - "EPBDS - tests extended: references on array"
#
# 1. C75: extend the IF() array formula with a third (false) argument,
#    which changes the computed sum from 7 to 9.
# 2. Add a new block of rows (97-99, 102-103) demonstrating an array
#    formula (SQRT of a 2x2 literal array) that spills into A97:B98, a
#    normal formula referencing a spilled array cell (C99), and a second
#    array formula that itself references the first array's range
#    (B102:C103 = A97:B98+1), together with the mirrored "I/J" value
#    columns and "N" labels used throughout this sheet.
# 3. Move the trailing "END" marker row from row 101 to row 110, and add a
#    second "END" label next to it (H110 / I110).
# 4. Give column C a custom width and move the visible selection down to
#    the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Update the existing array formula in C75 -------------------------
$ws.Range("C75").FormulaArray = "=SUM(IF(A67:B68>2,A67:B68,1))"

# --- 2. New "reference on array" example ---------------------------------
# A97:B98 <- array formula SQRT({1,2;3,4})
$ws.Range("A97:B98").FormulaArray = "=SQRT({1,2;3,4})"

# C99 <- ordinary formula that uses a cell from the spilled array (B98)
$ws.Range("C99").Formula = "=B98+1"
$ws.Range("I99").Value = 3
$ws.Range("N99").Value = "Use array cell"

# B102:C103 <- array formula referencing the A97:B98 array range
$ws.Range("B102:C103").FormulaArray = "=A97:B98+1"
$ws.Range("I102").Value = 2
$ws.Range("J102").Value = 2.4142135623730949
$ws.Range("N102").Value = "ref array in array"
$ws.Range("I103").Value = 2.7320508075688772
$ws.Range("J103").Value = 3

# --- 3. Move the "END" marker from H101 down to row 110 ------------------
$ws.Range("H101").ClearContents()
$ws.Range("H110").Value = "END"
$ws.Range("I110").Value = "END"

# --- 4. Cosmetics: column width + scrolled selection ----------------------
# (13.140625 is the stored/XML column width; ColumnWidth is specified in
# characters and gets snapped to the nearest whole pixel by Excel, so 12.25
# is the input that lands closest to the target stored width.)
$ws.Columns.Item(3).ColumnWidth = 12.25
$ws.Range("J103").Select()
